$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 08:57:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 08:57:28"
$wsZhCn.Range("K3").Value = "2016-09-01 08:57:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 08:57:32"
$wsDeDe.Range("K3").Value = "2016-09-01 08:57:54"
